$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card16")

# Header row: M1 loses its trailing space, N1 is the new "Correction " header
$ws.Cells.Item(1, 13).Value = "Event"
$ws.Cells.Item(1, 14).Value = "Correction "

# Give N1 the exact same formatting (bold, centered, bordered header style)
# as the rest of row 1 by copying formats from M1 -> N1.
$ws.Cells.Item(1, 13).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-12: M gets the literal "nan" text (previously blank),
# and a brand-new, still-blank N cell appears alongside it.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"

    # Force the new N cell into existence as an empty text cell (mirrors
    # the blank inlineStr cell M used to have) using the quote-prefix
    # trick, then strip the quote-prefix style back off again.
    $ws.Cells.Item($r, 14).Value = "'"
    $ws.Cells.Item($r, 14).Style = "Normal"
}
